# Update column G ("K") values on Sheet1 to reflect regenerated
# save_data (K computed instead of Strike#, with refreshed std/mean
# derived s_vals written back into the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$updates = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 4
    6  = 1
    7  = 3
    8  = 3
    9  = 3
    10 = 4
    12 = 1
    14 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
